$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the observation text on row 8 (G8)
$ws.Range("G8").Value = "Consumió 500 gigas de ram. Dejó 1896 vbles"

# Mark row 9 (the "Training strategy" step) as done
$ws.Range("A9").Value = "SI"

# Insert a new blank row after the "BO" row (row 10), for the new "Modelo final" step
$ws.Rows("11").Insert()

# Mark row 10 (the "BO" step) as done
$ws.Range("A10").Value = "SI"

# New row 11 with the final model note
$ws.Range("B11").Value = "Modelo final"

# Note the duration of the "BO" step
$ws.Range("G10").Value = "23 horas"

# Update selection to reflect the last edited cell
$ws.Range("G10").Select()
